$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
if ($ws -eq $null) {
    $ws = $wb.Worksheets.Item("cancer_de_mama")
}

# The 2025 data row (row 12) was removed from the "cancer_de_mama" sheet.
# Deleting the whole row shifts every row below it up by one, which also
# updates the sheet dimension and the shared-formula ranges that spanned it.
$ws.Rows.Item(12).Delete()

# Leave the selection where Excel lands after a row delete: on the row
# that slid up into the now-last (blank) row position.
[void]$ws.Range("A12:XFD12").Select()
